# Add an INSTRUCTORS column (I) to the "NewSubject" sheet with per-row
# instructor-load values, and move the "selected tab" focus from Student
# (sheet4) to NewSubject (sheet7), matching the upstream commit
# "add instructors to newsubject xslx".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewSubject")

# New header in column I (pushes sharedStrings uniqueCount 75 -> 76).
$ws.Range("I1").Value = "INSTRUCTORS"

# Per-row instructor-load figures for rows 2-9.
$ws.Range("I2").Value = 0.1
$ws.Range("I3").Value = 2
$ws.Range("I4").Value = 3.4
$ws.Range("I5").Value = 3.4
$ws.Range("I6").Value = 5.6
$ws.Range("I7").Value = 5
$ws.Range("I8").Value = 7
$ws.Range("I9").Value = 8

# The workbook now opens on the NewSubject tab (was Student before), with
# the cursor left on L8 there.
$ws.Activate()
$ws.Range("L8").Select()
